$wb = $excel.ActiveWorkbook

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 1844.1428
$ws.Cells.Item(100, 9).Value = 2261.8
$ws.Cells.Item(100, 11).Value = 2261.8
$ws.Cells.Item(100, 13).Value = -1720.8

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 50004724
$ws.Cells.Item(113, 9).Value = 66669970
$ws.Cells.Item(113, 10).Value = 8999
$ws.Cells.Item(113, 11).Value = 66669970
$ws.Cells.Item(113, 12).Value = 8999
$ws.Cells.Item(113, 13).Value = -66666716
$ws.Cells.Item(113, 14).Value = -15507

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 863
$ws.Cells.Item(129, 10).Value = 886.08887
$ws.Cells.Item(129, 12).Value = 2658.26661
$ws.Cells.Item(129, 14).Value = -12658.26661

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 7939049
$ws.Cells.Item(132, 9).Value = 9261627
$ws.Cells.Item(132, 10).Value = 3583.3333
$ws.Cells.Item(132, 11).Value = 27784881
$ws.Cells.Item(132, 12).Value = 10749.9999
$ws.Cells.Item(132, 13).Value = -27782351
$ws.Cells.Item(132, 14).Value = -15809.9999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1896.8235
$ws.Cells.Item(137, 9).Value = 1121.5555
$ws.Cells.Item(137, 10).Value = 2769
$ws.Cells.Item(137, 11).Value = 3364.6665
$ws.Cells.Item(137, 12).Value = 8307
$ws.Cells.Item(137, 13).Value = -814.6664999999998
$ws.Cells.Item(137, 14).Value = -13407

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1908.3877
$ws.Cells.Item(138, 9).Value = 756.64703
$ws.Cells.Item(138, 10).Value = 2150.111
$ws.Cells.Item(138, 11).Value = 2269.94109
$ws.Cells.Item(138, 12).Value = 6450.333
$ws.Cells.Item(138, 13).Value = 2870.05891
$ws.Cells.Item(138, 14).Value = -16730.333

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3498.7058
$ws.Cells.Item(32, 9).Value = 3417.449
$ws.Cells.Item(32, 11).Value = 3417.449
$ws.Cells.Item(32, 13).Value = -3130.449

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1903.6666
$ws.Cells.Item(61, 9).Value = 1440
$ws.Cells.Item(61, 10).Value = 2135.5
$ws.Cells.Item(61, 11).Value = 1440
$ws.Cells.Item(61, 12).Value = 2135.5
$ws.Cells.Item(61, 13).Value = -1228
$ws.Cells.Item(61, 14).Value = -2559.5

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2703
$ws.Cells.Item(74, 9).Value = 2506
$ws.Cells.Item(74, 10).Value = 2900
$ws.Cells.Item(74, 11).Value = 2506
$ws.Cells.Item(74, 12).Value = 2900
$ws.Cells.Item(74, 13).Value = -1632
$ws.Cells.Item(74, 14).Value = -4648

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2703
$ws.Cells.Item(77, 9).Value = 2506
$ws.Cells.Item(77, 10).Value = 2900
$ws.Cells.Item(77, 11).Value = 12530
$ws.Cells.Item(77, 12).Value = 14500
$ws.Cells.Item(77, 13).Value = -8162
$ws.Cells.Item(77, 14).Value = -23236

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1520.9333
$ws.Cells.Item(110, 9).Value = 1138.5385
$ws.Cells.Item(110, 11).Value = 1138.5385
$ws.Cells.Item(110, 13).Value = 906.4614999999999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1995.0588
$ws.Cells.Item(132, 9).Value = 1703.3864
$ws.Cells.Item(132, 11).Value = 5110.1592
$ws.Cells.Item(132, 13).Value = -2580.1592

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1903.6666
$ws.Cells.Item(136, 9).Value = 1440
$ws.Cells.Item(136, 10).Value = 2135.5
$ws.Cells.Item(136, 11).Value = 4320
$ws.Cells.Item(136, 12).Value = 6406.5
$ws.Cells.Item(136, 13).Value = -1770
$ws.Cells.Item(136, 14).Value = -11506.5

# BSM row 40
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(40, 8).Value = 35000
$ws.Cells.Item(40, 10).Value = 35000
$ws.Cells.Item(40, 12).Value = 35000
$ws.Cells.Item(40, 14).Value = -35530

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 16667487
$ws.Cells.Item(94, 9).Value = 22727918
$ws.Cells.Item(94, 11).Value = 22727918
$ws.Cells.Item(94, 13).Value = -22727467

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4993.2
$ws.Cells.Item(134, 9).Value = 1684.92
$ws.Cells.Item(134, 11).Value = 5054.76
$ws.Cells.Item(134, 13).Value = -2519.76

# CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 13).Value = $null

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 577.5
$ws.Cells.Item(122, 9).Value = 461.2
$ws.Cells.Item(122, 10).Value = 771.3333
$ws.Cells.Item(122, 11).Value = 1383.6
$ws.Cells.Item(122, 12).Value = 2313.9999
$ws.Cells.Item(122, 13).Value = 1066.4
$ws.Cells.Item(122, 14).Value = -7213.9999

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 5499.032
$ws.Cells.Item(132, 9).Value = 5985.75
$ws.Cells.Item(132, 11).Value = 17957.25
$ws.Cells.Item(132, 13).Value = -15427.25

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 11906198
$ws.Cells.Item(134, 9).Value = 14494180
$ws.Cells.Item(134, 10).Value = 1478.8
$ws.Cells.Item(134, 11).Value = 43482540
$ws.Cells.Item(134, 12).Value = 4436.4
$ws.Cells.Item(134, 13).Value = -43480005
$ws.Cells.Item(134, 14).Value = -9506.4

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 13).Value = $null

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 19609260
$ws.Cells.Item(131, 10).Value = 1511.9773
$ws.Cells.Item(131, 12).Value = 4535.9319
$ws.Cells.Item(131, 14).Value = -14615.9319

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 1962.8
$ws.Cells.Item(136, 9).Value = 999.2
$ws.Cells.Item(136, 10).Value = 2926.4
$ws.Cells.Item(136, 11).Value = 2997.6
$ws.Cells.Item(136, 12).Value = 8779.200000000001
$ws.Cells.Item(136, 13).Value = 2102.4
$ws.Cells.Item(136, 14).Value = -18979.2

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1791.5834
$ws.Cells.Item(61, 9).Value = 1214.2941
$ws.Cells.Item(61, 10).Value = 3193.5715
$ws.Cells.Item(61, 11).Value = 1214.2941
$ws.Cells.Item(61, 12).Value = 3193.5715
$ws.Cells.Item(61, 13).Value = -1012.2941
$ws.Cells.Item(61, 14).Value = -3597.5715

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1279.375
$ws.Cells.Item(68, 9).Value = 1022.6667
$ws.Cells.Item(68, 10).Value = 2049.5
$ws.Cells.Item(68, 11).Value = 1022.6667
$ws.Cells.Item(68, 12).Value = 2049.5
$ws.Cells.Item(68, 13).Value = -273.6667
$ws.Cells.Item(68, 14).Value = -3547.5

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 1279.375
$ws.Cells.Item(71, 9).Value = 1022.6667
$ws.Cells.Item(71, 10).Value = 2049.5
$ws.Cells.Item(71, 11).Value = 5113.3335
$ws.Cells.Item(71, 12).Value = 10247.5
$ws.Cells.Item(71, 13).Value = -1369.3335
$ws.Cells.Item(71, 14).Value = -17735.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1791.5834
$ws.Cells.Item(113, 9).Value = 1214.2941
$ws.Cells.Item(113, 10).Value = 3193.5715
$ws.Cells.Item(113, 11).Value = 1214.2941
$ws.Cells.Item(113, 12).Value = 3193.5715
$ws.Cells.Item(113, 13).Value = 955.7058999999999
$ws.Cells.Item(113, 14).Value = -7533.5715

# LTW row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 47429
$ws.Cells.Item(123, 10).Value = 47429
$ws.Cells.Item(123, 12).Value = 47429
$ws.Cells.Item(123, 14).Value = -57229

# WVR row 102
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).Value = $null

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 630.625
$ws.Cells.Item(136, 9).Value = 531.2353000000001
$ws.Cells.Item(136, 10).Value = 872
$ws.Cells.Item(136, 11).Value = 1593.7059
$ws.Cells.Item(136, 12).Value = 2616
$ws.Cells.Item(136, 13).Value = 956.2940999999998
$ws.Cells.Item(136, 14).Value = -7716
